# daily auto push: 2026-01-29 14:04 UTC
# Two new rows of data for 2026/01/29 (Thursday / 木) are appended to the
# existing block of same-date rows (rows 741-744), which pushes all of the
# subsequent rows (previously 745-786) down by two rows (now 747-788).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 745, shifting rows
# 745:786 down to 747:788.
$ws.Rows("745:746").Insert()

# Seed column A (date text) by copying an existing "2026/01/29" text cell
# so the new cells stay plain text instead of being auto-converted into a
# date serial number / receiving a new number-format style.
$ws.Range("A744").Copy($ws.Range("A745"))
$ws.Range("A744").Copy($ws.Range("A746"))

# Row 745: 2026/01/29, 木, 18:00, ranking 171
$ws.Range("B745").Value = "木"
$ws.Range("C745").Value = 18
$ws.Range("D745").Value = 171

# Row 746: 2026/01/29, 木, 22:00, ranking 195
$ws.Range("B746").Value = "木"
$ws.Range("C746").Value = 22
$ws.Range("D746").Value = 195
